$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.258.49"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "1.594.69"
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'211.58"
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  +0.30%  "
$ws.Range("D10").Value = "'18.97"
$ws.Range("E10").Value = "  -1.07%  "
$ws.Range("D11").Value = "'0.0853"
$ws.Range("E11").Value = "  +0.79%  "
$ws.Range("D12").Value = "1.819.14"
$ws.Range("D13").Value = "1.613.42"
$ws.Range("E13").Value = "  +1.74%  "
$ws.Range("E14").Value = "  -0.06%  "
$ws.Range("E15").Value = "  -2.47%  "
$ws.Range("D16").Value = "'63.65"
$ws.Range("E16").Value = "  -0.34%  "
$ws.Range("D17").Value = "26.243.32"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("D18").Value = "'230.20"
$ws.Range("E18").Value = "  +7.68%  "
$ws.Range("D19").Value = "'7.68"
$ws.Range("E19").Value = "  +5.37%  "
$ws.Range("E20").Value = "  -0.28%  "
$ws.Range("D21").Value = "'0.999"
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("D23").Value = "'2.15"
$ws.Range("E23").Value = "  +1.17%  "
$ws.Range("E24").Value = "  -0.49%  "
$ws.Range("D25").Value = "'145.75"
$ws.Range("E25").Value = "  +1.07%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").Value = "'7.00"
$ws.Range("E27").Value = "  +0.38%  "
$ws.Range("E28").Value = "  +0.47%  "
$ws.Range("E29").Value = "  +1.68%  "
$ws.Range("E30").Value = "  -0.44%  "
$ws.Range("E31").Value = "  +0.35%  "
$ws.Range("E32").Value = "  +0.76%  "
$ws.Range("D33").Value = "1.463.86"
$ws.Range("E33").Value = "  +3.27%  "
$ws.Range("E34").Value = "  +0.37%  "
$ws.Range("E35").Value = "  -0.36%  "
$ws.Range("E36").Value = "  +0.77%  "
$ws.Range("D37").Value = "'0.568"
$ws.Range("E37").Value = "  -3.21%  "
$ws.Range("E38").Value = "  -1.14%  "
$ws.Range("D39").Value = "'0.821"
$ws.Range("E39").Value = "  +0.19%  "
$ws.Range("E40").Value = "  -2.21%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("E42").Value = "  +2.11%  "
$ws.Range("D43").Value = "'0.931"
$ws.Range("E43").Value = "  -1.88%  "
$ws.Range("D44").Value = "1.731.70"
$ws.Range("E44").Value = "  +0.68%  "
$ws.Range("D45").Value = "'0.756"
$ws.Range("E45").Value = "  -1.20%  "
$ws.Range("D46").Value = "'60.46"
$ws.Range("E46").Value = "  -0.69%  "
$ws.Range("E47").Value = "  +2.73%  "
$ws.Range("E48").Value = "  -0.45%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").Value = "'0.998"
$ws.Range("E50").Value = "  -0.12%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'7.42"
$ws.Range("E51").Value = "  +0.90%  "
